# Daily attendance processing - 2026-01-06 01:38:43
# Swap the order of the two "Recorded By" contributors (column G) for the
# rows where the recorder list toggled between "System, dnasr281@gmail.com"
# and "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value changes from
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
$rowsToDnasrFirst = @(2,3,4,5,6,7,16,17,22,23,37,38,43,44,58,59,64,65,79,80,85,86,87,88,89,90,99,100,105,106,107,108,109,110,119,120,125,126,127,128,129,130,139,140,145,146,147,148,149,150,159,160,165,166,167,168,169,170,179,180,185,186,200,201,206,207,221,222,227,228,242,243)

# Rows whose "Recorded By" value changes from
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
$rowsToSystemFirst = @(8,24,26,29,45,47,50,66,68,71,91,111,131,151,171,187,189,192,208,210,213,229,231,234)

foreach ($r in $rowsToDnasrFirst) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

foreach ($r in $rowsToSystemFirst) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}
